$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at row 295, pushing the existing rows 295-315
# (and their D/J/K/L/M/O/P values) down to 296-316.
$ws.Rows.Item(295).EntireRow.Insert()

# Populate the newly inserted row 295 with the new weekly record.
$ws.Range("A295").Value = 10
$ws.Range("B295").Value = "Vega Modelo de Temuco"
$ws.Range("C295").Value = "La Araucanía"
$ws.Range("D295").Value = 44714
$ws.Range("E295").Value = 9
$ws.Range("F295").Value = 100112044
$ws.Range("G295").Value = "Perejil"
$ws.Range("H295").Value = "Sin especificar"
$ws.Range("I295").Value = "Primera"
$ws.Range("J295").Value = 100
$ws.Range("K295").Value = 4000
$ws.Range("L295").Value = 6000
$ws.Range("M295").Value = 5300
$ws.Range("N295").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O295").Value = "Provincia de Cautín"
$ws.Range("P295").Value = 1767
$ws.Range("Q295").Value = 3
$ws.Range("R295").Value = "Hortaliza"
